$wb = $excel.ActiveWorkbook

# "2016" sheet: update January purchase-price (D) and December dividend (R)
# figures for three holdings; the dependent yield/total formulas (E, S, T)
# and the row-12 sums recalc automatically. The "2017" sheet pulls its G
# column from '2016'!S, so its E/U/V cells ripple through too.
$ws2016 = $wb.Worksheets.Item("2016")
$ws2016.Activate()

$ws2016.Range("D4").Value = 20.352
$ws2016.Range("R4").Value = 9.69

$ws2016.Range("D7").Value = 10.07
$ws2016.Range("R7").Value = 5.96

$ws2016.Range("D8").Value = 8.124
$ws2016.Range("R8").Value = 7.54

# Matches the saved cursor position recorded in the sheet's <selection>.
$ws2016.Range("I15").Select()
